$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.987.12'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.70%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.952.50'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.87%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.53'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.74'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.51%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.949.34'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.79%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.22%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.27'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +4.02%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +6.76%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.96%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000236'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +5.33%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.73%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.66%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.439.67'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.84%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.924.35'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.56%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.72'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.49%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.946.19'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.85%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '442.27'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.55%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.46'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.04%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.00%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.29%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.26'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.80%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.12'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.78%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.14'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.87%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.82'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.02%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = 'NEARProtocol'
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.25'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +5.05%  '
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.20'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.68%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.48%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +16.63%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.42'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.73%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.993'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.60%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +6.16%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.60'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.31%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.77%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '49.73'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.27%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.52'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.50%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.117'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -4.28%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.281'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.21%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '39.01'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -7.45%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '135.53'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.42%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.694.44'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.18%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.15%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '362.15'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.27%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.48%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '22.78'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -3.35%  '
